$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26: AI를 활용한 고용률 예측 모델 개발기(1) -> (2), and link update
$ws.Range("D26").Value = "AI를 활용한 고용률 예측 모델 개발기(2)"
$ws.Range("E26").Value = "https://blog.est.ai/2021/04/employment-rate-2/"

# Row 32: Edit Distance (Levenshtein Distance) (퍼옴) -> 샘플링 (Sampling), and link update
$ws.Range("D32").Value = "샘플링 (Sampling)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/311"

# Row 39: Using tf.Print() in TensorFlow -> Top 7 Feature Selection Techniques in Machine Learning, and link update
$ws.Range("D39").Value = "Top 7 Feature Selection Techniques in Machine Learning"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Top-7-Feature-Selection-Techniques-in-Machine-Learning-1"

# Row 51: [python+openpose] ... -> [html] 이메일 주소에 링크 걸 때, mailto:, and link update
$ws.Range("D51").Value = "[html] 이메일 주소에 링크 걸 때, mailto:"
$ws.Range("E51").Value = "https://bskyvision.com/950"

$wb.Save()
